# Update the "Förändrad" (Changed) date column (C) from 45184 to 45186 for every
# data row, and append the "Beteckning" (column A) text as the friendly-name
# second argument of every HYPERLINK(...) formula in columns S, T, V, W, X, Y
# that does not already have one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Row + $usedRange.Rows.Count - 1

$hyperlinkCols = @("S", "T", "V", "W", "X", "Y")

for ($r = 2; $r -le $lastRow; $r++) {

    # --- Column C: bump the changed-date serial from 45184 to 45186 ---
    $cCell = $ws.Range("C" + $r)
    $cVal = $cCell.Value2
    if ($cVal -eq 45184) {
        $cCell.Value = 45186
    }

    # --- Beteckning (used as the HYPERLINK friendly name) ---
    $beteckning = $ws.Range("A" + $r).Value2

    foreach ($col in $hyperlinkCols) {
        $cell = $ws.Range($col + $r)
        $f = $cell.Formula
        if ($f -and $f.ToUpper().StartsWith("=HYPERLINK(") -and -not $f.Contains(",")) {
            # Strip the trailing ")" and append the second argument.
            $newFormula = $f.Substring(0, $f.Length - 1) + ', "' + $beteckning + '")'
            $cell.Formula = $newFormula
        }
    }
}
